$wb = $excel.ActiveWorkbook

# --- Update the existing " listOfGames Opt2" sheet ---
$ws3 = $wb.Worksheets.Item(" listOfGames Opt2")
$ws3.Range("E38").Select()

# --- Create the new "get_participants Opt1" sheet as a copy of " listOfGames Opt2" ---
$ws3.Copy($null, $ws3)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "get_participants Opt1"

# Update values on the new sheet (generator-based participant timings)
$newSheet.Range("B2").Value = 28.92
$newSheet.Range("C2").Value = 28.48
$newSheet.Range("D2").Value = 28.481000000000002

$newSheet.Range("B3").Value = 27.951000000000001
$newSheet.Range("C3").Value = 27.878
$newSheet.Range("D3").Value = 27.512

$newSheet.Range("B4").Value = 0.96199999999999997
$newSheet.Range("C4").Value = 0.95399999999999996
$newSheet.Range("D4").Value = 0.95899999999999996

$newSheet.Activate()
$newSheet.Range("B5").Select()

# Add the note about timings to the original sheet
$ws3.Range("B9").Value = "Note: This was after first making the changes. For some reason, the next morning all the timings are coming out around 33 seconds.."

# The original sheet is no longer the active tab, so (like the other inactive
# sheets in this workbook) it now carries an explicit pageSetup element.
$ws3.PageSetup.Orientation = 1
